# Automatische test-sync: 2025-07-22 12:55:50
# Appends the new test-mail log entry (#20) to the "Logs" sheet and
# refreshes the "Dashboard" category summary to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append row 19 with the new mail record
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A19").Value = "Kun je dit bespreken met finance?"
$logs.Range("B19").Value = "mailmind.test@zohomail.eu"
$logs.Range("C19").Value = "Testmail #20: Kun je dit bespreken met finance?"
$logs.Range("D19").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E19").Value = "Beste afzender,`nDank voor je e-mail. Om je beter van dienst te kunnen zijn, zou je meer specifieke details willen geven over wat er precies besproken moet worden met de financiële afdeling? `nMet vriendelijke groet,`n[Naam] `nE-mailassistent"
$logs.Range("F19").Value = "2025-07-22 12:54:57"
$logs.Range("G19").Value = "Ja"
$logs.Range("H19").Value = "Nee"
$logs.Range("I19").Value = "Ja"
$logs.Range("J19").Value = "Ja"

# Re-fit the new row's height now that it holds a multi-line answer so it
# doesn't keep an explicit/custom row height (matches the rest of the sheet).
$logs.Rows.Item(19).AutoFit()

# Extend the conditional-formatting ranges (D/G/H/I/J) so the new row
# is covered the same way the rest of the table is.
$colsToExtend = "D", "G", "H", "I", "J"
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2:" + $col + "18")
    $newRange = $logs.Range($col + "2:" + $col + "19")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: re-sync the category counts now that "Intern
#    verzoek / Actie voor medewerker" has gone from 1 to 2 occurrences,
#    which re-sorts it above "Overig" (2) and "Bestelling / Levering" (1)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B4").Value = 2

$dash.Range("A5").Value = "Overig"
$dash.Range("B5").Value = 2

$dash.Range("A7").Value = "Bestelling / Levering"
$dash.Range("B7").Value = 1
